$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new "CH model" columns: Standard_CH (E) and Rounded_CH (F) ---
$ws.Range("E1").Value = "Standard_CH"
$ws.Range("F1").Value = "Rounded_CH"

# Give the new headers the same look (bold / bordered / centered) as the
# existing header cells by copying the formatting from D1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Recomputed "Nash" (C) values plus the new Standard_CH (E) / Rounded_CH (F) values ---
# Row layout: row, C (Nash), E (Standard_CH), F (Rounded_CH)
$rows = @(
    @(2,  6.169729885736025,  4.334391407037445, 4.334391407037446),
    @(3,  5.993939155417372,  4.334391407037445, 4.334391407037446),
    @(4,  4.193396040772955,  4.334391407037445, 4.334391407037446),
    @(5,  2.1331024440159,    4.334391407037445, 4.334391407037446),
    @(6,  16.38592624660553, 16.71836685571586, 10.99047733666301),
    @(7,  7.019822372650111,  4.334391407037445, 4.334391407037446),
    @(8,  5.485589095727236,  4.334391407037445, 4.334391407037446),
    @(9,  5.869311398848177, 5.262587817610159, 10.99047733666301),
    @(10, 3.686357852204599,  4.334391407037445, 4.334391407037446),
    @(11, 7.845645243811248,  4.334391407037445, 4.334391407037446),
    @(12, 4.225536883256589,  4.334391407037445, 4.334391407037446),
    @(13, 6.838784157844334,  4.334391407037445, 4.334391407037446),
    @(14, 3.223923281776798,  4.334391407037445, 4.334391407037446),
    @(15, 4.51677227229197,   4.334391407037445, 4.334391407037446),
    @(16, 2.098993794849594,  4.334391407037445, 4.334391407037446),
    @(17, 2.631613470292678,  4.334391407037445, 4.334391407037446),
    @(18, 4.318023797342149,  4.334391407037445, 4.334391407037446),
    @(19, 2.405971637346347,  4.334391407037445, 4.334391407037446),
    @(20, 2.583074238786782,  4.334391407037445, 4.334391407037446),
    @(21, 2.374486730423603,  4.334391407037445, 4.334391407037446)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 3).Value = $r[1]   # C -> Nash
    $ws.Cells.Item($rowNum, 5).Value = $r[2]   # E -> Standard_CH
    $ws.Cells.Item($rowNum, 6).Value = $r[3]   # F -> Rounded_CH
}
